# Update "想去人数" (want-to-go count) figures for the two sheets that
# carry the full data table: "展览" and "全部类型".
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 1071
    $ws.Range("F3").Value = 30
    $ws.Range("F4").Value = 506
}
